$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MaccorTimeFunction")
$ws.Range("A1").Value = "Test"
